$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.92%"
$ws.Range("E2").Style = "Normal"
$ws.Range("E3").Value = "'-0.24%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.043"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.73%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08053"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-0.55%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.863"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-4.12%"
$ws.Range("E6").Style = "Normal"
$ws.Range("B7").Value = "'KuCoinToken"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "'7.779"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.55%"
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = "'MXToken"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'0.9252"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-1.61%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'0.1282"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-6.24%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "'WazirX"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'0.1898"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-0.32%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "'MandalaExchangeToken"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.09054"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-2.24%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "'BitrueCoin"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.03433"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-2.42%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "'BitMartToken"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'0.09864"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.57%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "'BitForexToken"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'0.001402"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-2.75%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "'TigerCash"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'0.006233"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'3.29%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'LEO"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'3.837"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'5.82%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'GateToken"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'4.117"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-1.34%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'12.51%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3416"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-0.41%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1302"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-3.30%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.813"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-7.65%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2409"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-5.12%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04364"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.99%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001228"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.79%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004857"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'2.15%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D27").Value = "'0.0001297"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-0.24%"
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'41.71%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D39").Value = "'0.01968"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-2.69%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05163"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-0.60%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007501"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-1.94%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.01013"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-9.41%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'-1.83%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002106"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'0.24%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009860"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-12.94%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006173"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-2.87%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.52%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'64.95"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-0.41%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.001247"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'4.43%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002096"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.52%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0001996"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.52%"
$ws.Range("E51").Style = "Normal"
